$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111, shifting existing rows 111-126 down to 112-127
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new data record
$ws.Range("A111").Value = 4
$ws.Range("B111").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C111").Value = "Los Lagos"
$ws.Range("D111").Value = 44491
$ws.Range("E111").Value = 10
$ws.Range("F111").Value = 100112039
$ws.Range("G111").Value = "Ciboulette"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 240
$ws.Range("K111").Value = 2500
$ws.Range("L111").Value = 2500
$ws.Range("M111").Value = 2500
$ws.Range("N111").Value = "$/docena de atados"
$ws.Range("O111").Value = "Región Metropolitana"
$ws.Range("P111").Value = 833
$ws.Range("Q111").Value = 3
$ws.Range("R111").Value = "Hortaliza"
